# framework setup with config, excel, report

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "Sheet1" to "Login"
$ws.Name = "Login"

# Fix the typo'd email address in A4 ("xy@abc.com" -> "xyz@abc.com")
$ws.Cells.Item(4, 1).Value = "xyz@abc.com"

# Move the active selection to B13
$ws.Range("B13").Select()
